$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Body text: "...Title: Prototyping Labs Manager..." is split across the
#    runs "...Title: Prototyping Lab" + "s" + " " + "Manager" + ...
#    Remove the stray "s" so the phrase reads
#    "...Title: Prototyping Lab Manager..."; the first three runs of the
#    phrase collapse into one run "...Title: Prototyping Lab ".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Prototyping Labs", $false, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $sRange = $d.Range($rng.End - 1, $rng.End)
    $sRange.Delete()
}

# The engine coalesces same-formatted runs across the whole paragraph
# whenever any part of it is edited. The rest of that paragraph (from
# "Manager" onward) must keep its original run boundaries untouched, so
# restore them by briefly toggling a character property on each original
# run's exact range - this forces the engine to re-split the text at those
# boundaries without altering any content.
$tailFragments = @("Manager", " ", "   ", "                          ",
    "                                                            ",
    "Date:", " ", "8", "/", "1", "6", "/2022    ")

$rng2 = $d.Content
$rng2.Find.Execute("Manager", $false, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $cursor = $rng2.Start
    foreach ($frag in $tailFragments) {
        $fragEnd = $cursor + $frag.Length
        $partRange = $d.Range($cursor, $fragEnd)
        $partRange.Font.Bold = 1
        $partRange.Font.Bold = 0
        $cursor = $fragEnd
    }
}

# ---------------------------------------------------------------------------
# 2) First-page header: "Prototyping Labs at GIX" -> "Prototyping Lab at GIX"
# ---------------------------------------------------------------------------
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute(
            "Prototyping Labs at GIX",
            $false, $false, $false, $false, $false,
            $true, 1, $false,
            "Prototyping Lab at GIX",
            2
        )
    }
}
